# Workbook/sheet handles.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FLUXOS")

# --- Add a new supplier row (row 3), mirroring the layout of row 2 ---
# COD FLUXO | COD IMS | COD FORNECEDOR | FORNECEDOR | COD DESTINO | NOME DESTINO
# | VEICULO PRINCIPAL | MOT | TECNOLOGIA | TIPO SATURACAO | PEDIDO | TRANSPORTADORA
$ws.Range("A3").Value = 76
$ws.Range("B3").Value = "32045/43513"
$ws.Range("C3").Value = 800043235
$ws.Range("D3").Value = "KOSTAL"
$ws.Range("E3").Value = 1080
$ws.Range("F3").Value = "FCA"
$ws.Range("G3").Value = "CARRETA"
$ws.Range("H3").Value = "FTL"
$ws.Range("I3").Value = "MONTAGEM-SP"
$ws.Range("J3").Value = "VOLUME"
$ws.Range("K3").Value = "D"
$ws.Range("L3").Value = "JAT"

# --- New working cell C4: empty, but pre-formatted bold + centered ---
$ws.Range("C4").Value = $null
$ws.Range("C4").HorizontalAlignment = -4108   # xlCenter
$ws.Range("C4").Font.Bold = $true
$ws.Range("C4").Font.Size = 11

# Leave the selection on the new working cell, scrolled into view.
$ws.Range("C4").Select()
try {
    $excel.ActiveWindow.ScrollColumn = 3
    $excel.ActiveWindow.ScrollRow = 1
} catch {
    # Best-effort only: view/scroll state isn't essential to the data edit.
}

$wb.Save()
